$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.138.24"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.835.29"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D5").Value = "242.10"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "0.6602"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "44.77"
$ws.Range("E8").Value = "  +6.31%  "
$ws.Range("D9").Value = "0.07396"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "23.11"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "0.07723"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.846.66"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "5.000"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "0.6706"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "82.26"
$ws.Range("E16").Value = "  -4.39%  "
$ws.Range("D17").Value = "6.144"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "0.000008679"
$ws.Range("E18").Value = "  +4.52%  "
$ws.Range("D19").Value = "29.150.23"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "2.091.99"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "225.10"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D24").Value = "7.148"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "1.002"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "158.37"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").Value = "8.582"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "0.1388"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").Value = "17.98"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "1.510"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "4.030"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "1.205"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "0.05388"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "1.848"
$ws.Range("D36").Value = "0.7449"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "1.297.96"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.761"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01792"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "6.354"
$ws.Range("E42").Value = "  +6.51%  "
$ws.Range("D43").Value = "0.9034"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "103.41"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("B47").Value = "XinFinNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D47").Value = "0.07923"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "1.990.53"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("D49").Value = "64.85"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "0.5137"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "1.747"
$ws.Range("E51").Value = "  -1.35%  "
